$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 4
$ws.Range("F3").Value = 2717
$ws.Range("F4").Value = 1055
$ws.Range("F5").Value = 19414
$ws.Range("F7").Value = 2197
$ws.Range("F9").Value = 611
$ws.Range("F10").Value = 423
$ws.Range("F11").Value = 681
$ws.Range("F12").Value = 232
$ws.Range("F13").Value = 247
$ws.Range("F15").Value = 359
$ws.Range("F16").Value = 71
$ws.Range("F17").Value = 257
$ws.Range("F19").Value = 182
$ws.Range("F22").Value = 92
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 189
$ws.Range("F7").Value = 280
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5989
$ws.Range("F3").Value = 635
$ws.Range("F4").Value = 583
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 5989
$ws.Range("F3").Value = 635
$ws.Range("F4").Value = 583
$ws.Range("F5").Value = 189
$ws.Range("F6").Value = 4
$ws.Range("F8").Value = 2717
$ws.Range("F9").Value = 1055
$ws.Range("F10").Value = 19414
$ws.Range("F15").Value = 280
$ws.Range("F16").Value = 2197
$ws.Range("F19").Value = 611
$ws.Range("F20").Value = 423
$ws.Range("F21").Value = 681
$ws.Range("F22").Value = 232
$ws.Range("F23").Value = 247
$ws.Range("F28").Value = 359
$ws.Range("F29").Value = 71
$ws.Range("F31").Value = 257
$ws.Range("F35").Value = 182
$ws.Range("F47").Value = 92
